$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row
# (row 2 .. row 381). Bump it from 46060 (2026-02-07) to 46061 (2026-02-08)
# for all of them in one shot.
$ws.Range("C2:C381").Value = 46061
